# Generate Report for Handoff
#
# dbed4300-f841-4f08-80b9-e7e2b6911b05 just got handed off (status moves
# from "In Translation" to "Ready for handoff" with a fresh handoff
# timestamp). The status report re-sorts the three affected rows
# (26abceb2..., d32f0acf..., dbed4300...) alphabetically by file name on
# every sheet, updating each row's hyperlink display text to match its
# new file name while leaving the underlying hyperlink targets and row
# styling untouched.

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value = "26abceb2-df30-4120-9bd6-8693b3848543.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "2016-03-25 00:52:31"
Set-LinkDisplay $ws '$A$7' "26abceb2-df30-4120-9bd6-8693b3848543.md"

$ws.Range("A8").Value = "d32f0acf-ec09-4e00-94d2-99390416cc32.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-03-25 00:48:41"
Set-LinkDisplay $ws '$A$8' "d32f0acf-ec09-4e00-94d2-99390416cc32.md"

$ws.Range("A9").Value = "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "2016-03-25 00:56:42"
Set-LinkDisplay $ws '$A$9' "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A7").Value = "26abceb2-df30-4120-9bd6-8693b3848543.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "26abceb2-df30-4120-9bd6-8693b3848543.9ca4aabde8470ac0fedcc05cb3302a82e2b36035.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-25 00:52:27"
Set-LinkDisplay $ws '$A$7' "26abceb2-df30-4120-9bd6-8693b3848543.md"
Set-LinkDisplay $ws '$D$7' "26abceb2-df30-4120-9bd6-8693b3848543.9ca4aabde8470ac0fedcc05cb3302a82e2b36035.zh-cn.xlf"

$ws.Range("A8").Value = "d32f0acf-ec09-4e00-94d2-99390416cc32.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "d32f0acf-ec09-4e00-94d2-99390416cc32.97ac89b367ddb70ca7e225c967c546e637c1671a.zh-cn.xlf"
$ws.Range("E8").Value = "2016-03-25 00:48:37"
Set-LinkDisplay $ws '$A$8' "d32f0acf-ec09-4e00-94d2-99390416cc32.md"
Set-LinkDisplay $ws '$D$8' "d32f0acf-ec09-4e00-94d2-99390416cc32.97ac89b367ddb70ca7e225c967c546e637c1671a.zh-cn.xlf"

$ws.Range("A9").Value = "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "dbed4300-f841-4f08-80b9-e7e2b6911b05.87ddbe91c5eabed32895a3ddbd3b4b8072f87e59.zh-cn.xlf"
$ws.Range("E9").Value = "2016-03-25 00:56:38"
Set-LinkDisplay $ws '$A$9' "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"
Set-LinkDisplay $ws '$D$9' "dbed4300-f841-4f08-80b9-e7e2b6911b05.87ddbe91c5eabed32895a3ddbd3b4b8072f87e59.zh-cn.xlf"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A7").Value = "26abceb2-df30-4120-9bd6-8693b3848543.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "26abceb2-df30-4120-9bd6-8693b3848543.9ca4aabde8470ac0fedcc05cb3302a82e2b36035.de-de.xlf"
$ws.Range("E7").Value = "2016-03-25 00:52:31"
Set-LinkDisplay $ws '$A$7' "26abceb2-df30-4120-9bd6-8693b3848543.md"
Set-LinkDisplay $ws '$D$7' "26abceb2-df30-4120-9bd6-8693b3848543.9ca4aabde8470ac0fedcc05cb3302a82e2b36035.de-de.xlf"

$ws.Range("A8").Value = "d32f0acf-ec09-4e00-94d2-99390416cc32.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "d32f0acf-ec09-4e00-94d2-99390416cc32.97ac89b367ddb70ca7e225c967c546e637c1671a.de-de.xlf"
$ws.Range("E8").Value = "2016-03-25 00:48:41"
Set-LinkDisplay $ws '$A$8' "d32f0acf-ec09-4e00-94d2-99390416cc32.md"
Set-LinkDisplay $ws '$D$8' "d32f0acf-ec09-4e00-94d2-99390416cc32.97ac89b367ddb70ca7e225c967c546e637c1671a.de-de.xlf"

$ws.Range("A9").Value = "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "dbed4300-f841-4f08-80b9-e7e2b6911b05.87ddbe91c5eabed32895a3ddbd3b4b8072f87e59.de-de.xlf"
$ws.Range("E9").Value = "2016-03-25 00:56:42"
Set-LinkDisplay $ws '$A$9' "dbed4300-f841-4f08-80b9-e7e2b6911b05.md"
Set-LinkDisplay $ws '$D$9' "dbed4300-f841-4f08-80b9-e7e2b6911b05.87ddbe91c5eabed32895a3ddbd3b4b8072f87e59.de-de.xlf"
